$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.363.01"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "1.667.67"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'312.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").Value = "'0.3946"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.79%  "
$ws.Range("D8").Value = "'0.3915"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("D9").Value = "'52.22"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.02%  "
$ws.Range("D10").Value = "'1.386"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.08%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "'0.08560"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").Value = "'24.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.33%  "
$ws.Range("D14").Value = "'7.275"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.93%  "
$ws.Range("D15").Value = "'7.983"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.18%  "
$ws.Range("D16").Value = "'0.00001331"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.68%  "
$ws.Range("D17").Value = "1.662.46"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("D18").Value = "'95.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "'0.07032"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("D20").Value = "'20.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").Value = "'6.985"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'13.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").Value = "24.373.85"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("D25").Value = "'2.515"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.97%  "
$ws.Range("D26").Value = "'3.076"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +14.20%  "
$ws.Range("D27").Value = "'22.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("D28").Value = "'156.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "'142.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("D30").Value = "'5.351"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Value = "'7.919"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.32%  "
$ws.Range("D32").Value = "'2.543"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.69%  "
$ws.Range("D33").Value = "1.846.92"
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("D34").Value = "'1.057"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.85%  "
$ws.Range("D35").Value = "'0.03112"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.25%  "
$ws.Range("D36").Value = "'0.08210"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").Value = "'6.860"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").Value = "'11.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.85%  "
$ws.Range("D39").Value = "'0.2758"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.48%  "
$ws.Range("D40").Value = "'0.09270"
$ws.Range("D40").Style = "Normal"
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'13.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.84%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.7659"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.70%  "
$ws.Range("D43").Value = "'1.445"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").Value = "'16.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("D45").Value = "'0.7072"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.97%  "
$ws.Range("D46").Value = "'2.529"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.95%  "
$ws.Range("D47").Value = "'4.125"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.11%  "
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").Value = "'0.08421"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").Value = "'136.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.03%  "
$ws.Range("D51").Value = "'1.258"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.14%  "
